$d = $word.ActiveDocument

# 1. Bold the "lf) Tomando el ejercicio de los deportistas, implementarlo todo con mvvm." sentence.
$r = $d.Content
$r.Find.Execute("lf) Tomando el ejercicio de los deportistas, implementarlo todo con mvvm.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Bold = 1
$r.BoldBi = 1

# 2. Move the "_GoBack" bookmark from its old location (end of "... en una grid. ")
#    to the new location inside "esos valores" -> "esos valore|s se debe mostrar...".
$r2 = $d.Content
$r2.Find.Execute("esos valore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $r2.End
$target = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $target)
